# Refactored exadata, added domain layer
$wb = $excel.ActiveWorkbook

# --- Storage sheet: drop the "[TB]" suffix from the Cell Raw / Total Used / Usable Available headers ---
$storage = $wb.Worksheets.Item("Storage")
$storage.Range("F1").Value = "Cell Raw"
$storage.Range("G1").Value = "Total Used"
$storage.Range("H1").Value = "Usable Available"

# --- Make "Storage" the active/selected sheet (instead of "Cluster View") ---
$storage.Activate()
$storage.Range("H1").Select()

# --- "Cluster View" keeps its own selection but is no longer the active tab ---
$clusterView = $wb.Worksheets.Item("Cluster View")
$clusterView.Range("E3").Select()

# re-activate Storage so it ends up as the workbook's active sheet/tab
$storage.Activate()
